$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 32
$ws.Range("I9").Value = 12.333333
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 12.333333
$ws.Range("L9").Value = 150
$ws.Range("M9").Value = 156.666667
$ws.Range("N9").Value = -488
$ws.Range("H28").Value = 234.09091
$ws.Range("I28").Value = 234.09091
$ws.Range("K28").Value = 234.09091
$ws.Range("M28").Value = 250.90909
$ws.Range("H41").Value = 219.83333
$ws.Range("I41").Value = 327.25
$ws.Range("K41").Value = 327.25
$ws.Range("M41").Value = 112.75
$ws.Range("H87").Value = 35000
$ws.Range("J87").Value = 35000
$ws.Range("L87").Value = 35000
$ws.Range("N87").Value = -37496
$ws.Range("H90").Value = 35000
$ws.Range("J90").Value = 35000
$ws.Range("L90").Value = 105000
$ws.Range("N90").Value = -117480
$ws.Range("H94").Value = 2033.3334
$ws.Range("J94").Value = 2500
$ws.Range("L94").Value = 2500
$ws.Range("N94").Value = -3402
$ws.Range("H116").Value = 5747.5
$ws.Range("I116").Value = 5747.5
$ws.Range("K116").Value = 5747.5
$ws.Range("M116").Value = -2305.5
$ws.Range("H132").Value = 3270.7
$ws.Range("I132").Value = 3270.7
$ws.Range("K132").Value = 9812.099999999999
$ws.Range("M132").Value = -7282.099999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1774.45
$ws.Range("I45").Value = 1771.25
$ws.Range("J45").Value = 1787.25
$ws.Range("K45").Value = 1771.25
$ws.Range("L45").Value = 1787.25
$ws.Range("M45").Value = -1394.25
$ws.Range("N45").Value = -2541.25
$ws.Range("H62").Value = 49999.5
$ws.Range("J62").Value = 49999.5
$ws.Range("L62").Value = 49999.5
$ws.Range("N62").Value = -51247.5
$ws.Range("H65").Value = 49999.5
$ws.Range("J65").Value = 49999.5
$ws.Range("L65").Value = 149998.5
$ws.Range("N65").Value = -156238.5
$ws.Range("H102").Value = 2532
$ws.Range("I102").Value = 2412.3333
$ws.Range("K102").Value = 2412.3333
$ws.Range("M102").Value = -790.3332999999998
$ws.Range("H132").Value = 987.1429000000001
$ws.Range("I132").Value = 987.1429000000001
$ws.Range("K132").Value = 2961.4287
$ws.Range("M132").Value = -431.4287000000004

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 5001.6665
$ws.Range("I5").Value = 2500
$ws.Range("J5").Value = 10005
$ws.Range("K5").Value = 2500
$ws.Range("L5").Value = 10005
$ws.Range("M5").Value = -2387
$ws.Range("N5").Value = -10231
$ws.Range("I20").Value = 2474.6667
$ws.Range("J20").Value = 5500
$ws.Range("K20").Value = 2474.6667
$ws.Range("L20").Value = 5500
$ws.Range("M20").Value = -2227.6667
$ws.Range("N20").Value = -5994
$ws.Range("H94").Value = 3589.3333
$ws.Range("I94").Value = 3589.3333
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 3589.3333
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -3138.3333
$ws.Range("N94").Value = ""
$ws.Range("H99").Value = 1069.4445
$ws.Range("I99").Value = 828.125
$ws.Range("K99").Value = 828.125
$ws.Range("M99").Value = 669.875
$ws.Range("H107").Value = 675.9
$ws.Range("I107").Value = 675.9
$ws.Range("K107").Value = 675.9
$ws.Range("M107").Value = 1244.1

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4459.9375
$ws.Range("I134").Value = 3335.3845
$ws.Range("K134").Value = 10006.1535
$ws.Range("M134").Value = -7471.1535
$ws.Range("H141").Value = 141649.83
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 141649.83
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 141649.83
$ws.Range("M141").Value = ""
$ws.Range("N141").Value = -152009.83

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1450.75
$ws.Range("J5").Value = 488.33334
$ws.Range("L5").Value = 1465.00002
$ws.Range("N5").Value = -1689.00002
$ws.Range("H23").Value = 422.66666
$ws.Range("J23").Value = 422.66666
$ws.Range("L23").Value = 1267.99998
$ws.Range("N23").Value = -1737.99998
$ws.Range("H59").Value = 100
$ws.Range("I59").Value = 100
$ws.Range("K59").Value = 300
$ws.Range("M59").Value = 240
$ws.Range("H135").Value = 1450.75
$ws.Range("J135").Value = 488.33334
$ws.Range("L135").Value = 4395.00006
$ws.Range("N135").Value = -9465.00006

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 38498
$ws.Range("J74").Value = 38498
$ws.Range("L74").Value = 38498
$ws.Range("N74").Value = -40370
$ws.Range("H77").Value = 38498
$ws.Range("J77").Value = 38498
$ws.Range("L77").Value = 115494
$ws.Range("N77").Value = -124854
$ws.Range("H102").Value = 1566
$ws.Range("I102").Value = 1566
$ws.Range("K102").Value = 1566
$ws.Range("M102").Value = 56
$ws.Range("H126").Value = 1333.6666
$ws.Range("I126").Value = 1000.5
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 3001.5
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -531.5
$ws.Range("N126").Value = -10940

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 874.6923
$ws.Range("I16").Value = 874.6923
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 874.6923
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -704.6923
$ws.Range("N16").Value = ""
$ws.Range("H62").Value = 49999.25
$ws.Range("J62").Value = 49999.25
$ws.Range("L62").Value = 49999.25
$ws.Range("N62").Value = -51247.25
$ws.Range("H65").Value = 49999.25
$ws.Range("J65").Value = 49999.25
$ws.Range("L65").Value = 149997.75
$ws.Range("N65").Value = -156237.75
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").Value = ""
$ws.Range("H122").Value = 3990
$ws.Range("I122").Value = 3486.6667
$ws.Range("K122").Value = 10460.0001
$ws.Range("M122").Value = -8010.000100000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 849.5
$ws.Range("I100").Value = 799.6667
$ws.Range("K100").Value = 1599.3334
$ws.Range("M100").Value = -1058.3334
$ws.Range("H109").Value = 119999
$ws.Range("J109").Value = 119999
$ws.Range("L109").Value = 119999
$ws.Range("N109").Value = -122773
